$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Every row in column A ("ConceptScheme") holds the shared label
#    "Data Modelling" (rows 2-47, one shared string). Rename it to
#    "DataModelling" (single word) everywhere it is used.
$ws.Range("A2:A47").Value = "DataModelling"

# 2. Flatten the rich-text "Datatype Scheme" definition (G30) down to a
#    single plain run - same wording, just without the split runs.
$ws.Range("G30").Value = "A Datatype scheme is a vocabulary of datatypes that might be used in a given context. Commonly, collections of datatype specifications will be grouped into a Datatype Scheme to form a recognisable standard. These could be proprietary to different standards bodies, database platforms, languages and data format specifications."

# 3. Shrink row 31 ("DataType") back down to its tighter auto height.
$ws.Rows.Item(31).RowHeight = 23.85

# 4. Scroll the sheet so column G is flush with the left edge, and move
#    the active selection to I7 (matches the saved view state).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 7
$ws.Range("I7").Select()
